$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (previously A2=3,B2=94 -> now A2=0,B2=228)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 228

# Update row 3 values (previously A3=0,B3=92 -> now A3=1,B3=87)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 87

# Remove rows 4 and 5 entirely (last two data rows dropped)
$ws.Range("A4:B5").Delete()
